# Update "想去人数" (F column) values on sheets "展览" and "全部类型"
# per the commit's re-generated data snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1252
    $ws.Range("F3").Value = 17114
    $ws.Range("F4").Value = 53
    $ws.Range("F5").Value = 1668
    $ws.Range("F6").Value = 76
    $ws.Range("F8").Value = 1048
    $ws.Range("F9").Value = 403
    $ws.Range("F11").Value = 136

    if ($sheetName -eq "展览") {
        $ws.Range("F12").Value = 11900
        $ws.Range("F14").Value = 52
        $ws.Range("F15").Value = 11594
        $ws.Range("F16").Value = 4716
        $ws.Range("F17").Value = 505
        $ws.Range("F19").Value = 415
        $ws.Range("F22").Value = 349
    }
    else {
        $ws.Range("F14").Value = 11900
        $ws.Range("F16").Value = 52
        $ws.Range("F17").Value = 11594
        $ws.Range("F18").Value = 4716
        $ws.Range("F19").Value = 505
        $ws.Range("F21").Value = 415
        $ws.Range("F24").Value = 349
    }
}
